# Update "想去人数" (number of people interested) figures for several
# events on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F6").Value = 268
    $ws.Range("F7").Value = 6494
    $ws.Range("F12").Value = 34
    $ws.Range("F16").Value = 517
}
